# DeveloperGuide: Rework UML diagrams
#
# The ":Address" lifeline object in the UndoRedo sequence diagram is
# renamed to ":Task" (the diagram was reworked to talk about the Task
# model object instead of Address). The shape is a two-paragraph
# rectangle ("Rectangle 62") whose first paragraph reads ":Address" and
# whose second paragraph reads "BookParser" (left untouched).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape = $s.Shapes.Item(5)

$tr = $shape.TextFrame.TextRange

# Only touch the first paragraph's characters (":Address" -> 8 chars)
# so the second paragraph ("BookParser") and its run/endParaRPr formatting
# stay exactly as they were.
$firstLine = $tr.Characters(1, 8)
$firstLine.Text = ":Task"
